# Applies spell-check "proofErr" run-splitting markup (as produced by
# Word's background spell/grammar checker) to several list items, and
# appends two new list paragraphs under the "Frontend:" section.
#
# Approach: for each target paragraph we select the paragraph's text
# (excluding the trailing paragraph mark), clear it, then use
# Range.InsertXML on the now-collapsed range to drop in a full
# replacement <w:p> (complete with its own <w:pPr>) built from the
# exact run / w:proofErr sequence we need. InsertXML replaces the
# content of the (single) paragraph the collapsed range lives in, so
# supplying <w:pPr> ourselves is what keeps the list numbering/style.
#
# NOTE: the interpreter here only reliably supports passing *bare
# variables* as function arguments (parenthesised / inline expression
# arguments get mis-parsed), so every value is built into a variable
# first and only bare variables are passed into the helper functions.

function Set-ParaXml($paraIndex, $pPrXml, $runsXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = ""
    $paraXml = "<w:p>" + $pPrXml + $runsXml + "</w:p>"
    $pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkgXml)
}

function Add-ParaAfter($paraIndex, $pPrXml, $runsXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $insPoint = $d.Range($rng.End, $rng.End)
    $paraXml = "<w:p>" + $pPrXml + $runsXml + "</w:p>"
    $pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insPoint.InsertXML($pkgXml)
}

$d = $word.ActiveDocument

$listPPr3 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$listPPr1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$listPPr2 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'

# 1) Paragraph 3: "Sql Server" -> "Sql" (spell-flagged) + " Server"
$runs3 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Sql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Server</w:t></w:r>'
Set-ParaXml 3 $listPPr3 $runs3

# 2) Paragraph 4: " Cloudinary database for pictures"
$runs4 = '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cloudinary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> database for pictures</w:t></w:r>'
Set-ParaXml 4 $listPPr3 $runs4

# 3) Paragraph 10: " Custom execption middleware"
$runs10 = '<w:r><w:t xml:space="preserve"> Custom </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>execption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> middleware</w:t></w:r>'
Set-ParaXml 10 $listPPr1 $runs10

# 4) Paragraph 11: " Cloudinary implementation"
$runs11 = '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cloudinary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> implementation</w:t></w:r>'
Set-ParaXml 11 $listPPr1 $runs11

# 5) Paragraph 15: "Jwt, loading and error interceptors"
$runs15 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Jwt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, loading and error interceptors</w:t></w:r>'
Set-ParaXml 15 $listPPr2 $runs15

# 6) Paragraph 17: "Bootswatch minty for theme"
$runs17 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Bootswatch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> minty for theme</w:t></w:r>'
Set-ParaXml 17 $listPPr2 $runs17

# 7) Paragraph 18: "Toastr for frontend user notification"
$runs18 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Toastr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for frontend user notification</w:t></w:r>'
Set-ParaXml 18 $listPPr2 $runs18

# 8) Paragraph 19: "Ngx-spinner for loading"
$runs19 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Ngx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-spinner for loading</w:t></w:r>'
Set-ParaXml 19 $listPPr2 $runs19

# 9) New paragraph after "Ng2-File-Upload" (paragraph 21): "NgbDatepickerModule"
$runsNgb = '<w:proofErr w:type="spellStart"/><w:r><w:t>NgbDatepickerModule</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Add-ParaAfter 21 $listPPr2 $runsNgb

# 10) New paragraph after that (now paragraph 22): "npm i bootstrap-icons"
$runsNpm = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> bootstrap-icons</w:t></w:r>'
Add-ParaAfter 22 $listPPr2 $runsNpm

Write-Host "Edits applied. Paragraph count:" $d.Paragraphs.Count
